# Populate the insurance ("保險") sheet (sheet5) with the full metadata
# columns, matching the pattern already used on the other sheets
# (land/building/deposit/stock): category, date, legislator_name,
# legislator_id, source_file, index.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

# --- Extend formatting to the new columns (F:K) before writing data ---
# Header row (B1) carries the bold/bordered header style -> copy it onto
# the new header cells.
$ws.Range("B1").Copy()
$ws.Range("F1:K1").PasteSpecial(-4122)

# Data rows (D2:D4) carry the plain data style -> copy it onto the new
# data cells.
$ws.Range("D2:D4").Copy()
$ws.Range("F2:K4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Header row ---
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# --- Row 2 (index 101) ---
$ws.Range("B2").Value = "國泰人壽"
$ws.Range("C2").Value = "鍾愛一生313"
$ws.Range("D2").Value = "黃靜秋"
$ws.Range("E2").Value = "insurance"
$ws.Range("F2").Value = "normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2011-12-31"
$ws.Range("H2").Value = "羅明才"
$ws.Range("I2").Value = 879
$ws.Range("J2").Value = "tmp94331"
$ws.Range("K2").Value = 101

# --- Row 3 (index 102) ---
$ws.Range("B3").Value = "保德信國際人壽"
$ws.Range("C3").Value = "教肓终身壽險"
$ws.Range("D3").Value = "黃靜秋"
$ws.Range("E3").Value = "insurance"
$ws.Range("F3").Value = "normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "2011-12-31"
$ws.Range("H3").Value = "羅明才"
$ws.Range("I3").Value = 879
$ws.Range("J3").Value = "tmp94331"
$ws.Range("K3").Value = 102

# --- Row 4 (index 103) ---
$ws.Range("B4").Value = "保德信國際人壽"
$ws.Range("C4").Value = "教育终身壽險"
$ws.Range("D4").Value = "黃靜秋"
$ws.Range("E4").Value = "insurance"
$ws.Range("F4").Value = "normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "2011-12-31"
$ws.Range("H4").Value = "羅明才"
$ws.Range("I4").Value = 879
$ws.Range("J4").Value = "tmp94331"
$ws.Range("K4").Value = 103

# Restore the plain data style on the date column cells now that the
# text value is safely stored (NumberFormat="@" above only existed to
# stop "2011-12-31" being auto-parsed into a date serial number).
$ws.Range("D2:D4").Copy()
$ws.Range("G2:G4").PasteSpecial(-4122)
$excel.CutCopyMode = 0
